# Update "F" column (想去人数 / interest count) values on three worksheets
# (展览, 演出, 全部类型) to reflect freshly re-generated counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 47
$ws1.Cells.Item(5, 6).Value = 181
$ws1.Cells.Item(6, 6).Value = 1073
$ws1.Cells.Item(8, 6).Value = 8142
$ws1.Cells.Item(9, 6).Value = 137
$ws1.Cells.Item(10, 6).Value = 205
$ws1.Cells.Item(11, 6).Value = 6877
$ws1.Cells.Item(14, 6).Value = 4972
$ws1.Cells.Item(15, 6).Value = 10
$ws1.Cells.Item(16, 6).Value = 5397
$ws1.Cells.Item(18, 6).Value = 329
$ws1.Cells.Item(19, 6).Value = 335
$ws1.Cells.Item(20, 6).Value = 459
$ws1.Cells.Item(26, 6).Value = 9132
$ws1.Cells.Item(28, 6).Value = 1652
$ws1.Cells.Item(29, 6).Value = 697
$ws1.Cells.Item(30, 6).Value = 41
$ws1.Cells.Item(32, 6).Value = 1853
$ws1.Cells.Item(33, 6).Value = 72
$ws1.Cells.Item(36, 6).Value = 1008
$ws1.Cells.Item(37, 6).Value = 1871
$ws1.Cells.Item(38, 6).Value = 241
$ws1.Cells.Item(40, 6).Value = 4762
$ws1.Cells.Item(41, 6).Value = 374
$ws1.Cells.Item(43, 6).Value = 74
$ws1.Cells.Item(48, 6).Value = 1248
$ws1.Cells.Item(50, 6).Value = 63

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(14, 6).Value = 88
$ws2.Cells.Item(17, 6).Value = 890

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 47
$ws4.Cells.Item(6, 6).Value = 181
$ws4.Cells.Item(8, 6).Value = 1073
$ws4.Cells.Item(10, 6).Value = 8143
$ws4.Cells.Item(11, 6).Value = 137
$ws4.Cells.Item(12, 6).Value = 205
$ws4.Cells.Item(13, 6).Value = 6877
$ws4.Cells.Item(17, 6).Value = 4972
$ws4.Cells.Item(18, 6).Value = 10
$ws4.Cells.Item(19, 6).Value = 5398
$ws4.Cells.Item(21, 6).Value = 329
$ws4.Cells.Item(22, 6).Value = 335
$ws4.Cells.Item(23, 6).Value = 459
$ws4.Cells.Item(27, 6).Value = 9132
$ws4.Cells.Item(29, 6).Value = 1652
$ws4.Cells.Item(30, 6).Value = 698
$ws4.Cells.Item(31, 6).Value = 41
$ws4.Cells.Item(33, 6).Value = 1853
$ws4.Cells.Item(34, 6).Value = 72
$ws4.Cells.Item(37, 6).Value = 1008
$ws4.Cells.Item(38, 6).Value = 1871
$ws4.Cells.Item(39, 6).Value = 241
$ws4.Cells.Item(41, 6).Value = 4762
$ws4.Cells.Item(42, 6).Value = 374
$ws4.Cells.Item(44, 6).Value = 74
$ws4.Cells.Item(48, 6).Value = 1248
$ws4.Cells.Item(50, 6).Value = 63
